# Landscaping Data.xlsx update:
# Append 7 new observation rows (534-540, all dated 45863 / 2025-07-25) to Sheet1,
# extending the data table that previously ended at row 533.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the date-cell number format (style index used by column A) down
#     onto the new rows before we set any values, so A534:A540 pick up the
#     same short-date display style as the rest of column A. ---
$ws.Range("A533").Copy()
$ws.Range("A534:A540").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Columns A:E (Date, Plant_Type, Plant_Size, Low, High) ---
$dataAE = @(
  @(45863, "Flowering",    "Large",  73, 89),
  @(45863, "Nonflowering", "Medium", 73, 89),
  @(45863, "Nonflowering", "Small",  73, 89),
  @(45863, "Nonflowering", "Medium", 73, 89),
  @(45863, "Nonflowering", "Medium", 73, 89),
  @(45863, "Nonflowering", "Large",  73, 89),
  @(45863, "Tree",         "Medium", 73, 89)
)
$arrAE = New-Object 'object[,]' 7, 5
for ($i = 0; $i -lt 7; $i++) {
  for ($j = 0; $j -lt 5; $j++) {
    $arrAE[$i, $j] = $dataAE[$i][$j]
  }
}
$ws.Range("A534:E540").Value = $arrAE

# --- Column F (Temp_Diff) keeps the existing ABS(Low-High) formula pattern,
#     continuing the shared-formula fill that starts at F478. ---
$ws.Range("F534:F540").Formula = "=ABS(D534-E534)"

# --- Columns G:T (Rain, Growth, Pruned, Quadrant, Shade, UV, Humidity,
#     Dew_Point, Pressure, Wind_Gust, Cloud_Cover, Visibility, AQI, Pollen) ---
$dataGT = @(
  @(0.18, 0.2,  "Yes", 2, "Bright",  6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.2,  "Yes", 3, "Bright",  6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.15, "Yes", 3, "Neutral", 6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.1,  "Yes", 3, "Neutral", 6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.2,  "Yes", 3, "Bright",  6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.4,  "Yes", 4, "Bright",  6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28),
  @(0.18, 0.95, "Yes", 1, "Neutral", 6, 0.82, 75, 30.07, 16, 0.79, 9.3, 57, 28)
)
$arrGT = New-Object 'object[,]' 7, 14
for ($i = 0; $i -lt 7; $i++) {
  for ($j = 0; $j -lt 14; $j++) {
    $arrGT[$i, $j] = $dataGT[$i][$j]
  }
}
$ws.Range("G534:T540").Value = $arrGT

# --- Update the sheet view: active cell / selection moves just past the new
#     data (I541), and the viewport scrolls so row 516 is at the top. ---
$win = $excel.ActiveWindow
$ws.Range("I541").Select()
$win.ScrollRow = 516
$win.ScrollColumn = 1
